$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $origStyle = $ws.Range($addr).Style
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = $origStyle
}

Set-TextValue $ws 'D2' '27.959.08'
Set-TextValue $ws 'E2' '  -2.25%  '
Set-TextValue $ws 'D3' '1.880.74'
Set-TextValue $ws 'E3' '  -1.67%  '
Set-TextValue $ws 'D4' '1.002'
Set-TextValue $ws 'E4' '  +0.15%  '
Set-TextValue $ws 'D5' '313.04'
Set-TextValue $ws 'E5' '  -0.60%  '
Set-TextValue $ws 'E6' '  +0.13%  '
Set-TextValue $ws 'D7' '0.5001'
Set-TextValue $ws 'E7' '  -4.12%  '
Set-TextValue $ws 'D8' '0.3846'
Set-TextValue $ws 'E8' '  -2.75%  '
Set-TextValue $ws 'D9' '0.09148'
Set-TextValue $ws 'E9' '  -5.90%  '
Set-TextValue $ws 'D10' '1.120'
Set-TextValue $ws 'E10' '  -2.97%  '
Set-TextValue $ws 'D11' '41.62'
Set-TextValue $ws 'E11' '  -1.12%  '
Set-TextValue $ws 'D12' '6.318'
Set-TextValue $ws 'E12' '  -3.50%  '
Set-TextValue $ws 'B13' 'WrappedEther'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws 'D13' '1.892.61'
Set-TextValue $ws 'E13' '  -1.12%  '
Set-TextValue $ws 'B14' 'Solana'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws 'D14' '20.70'
Set-TextValue $ws 'E14' '  -2.45%  '
Set-TextValue $ws 'D15' '7.255'
Set-TextValue $ws 'E15' '  -3.90%  '
Set-TextValue $ws 'E16' '  +0.20%  '
Set-TextValue $ws 'E17' '  -2.76%  '
Set-TextValue $ws 'D18' '91.44'
Set-TextValue $ws 'E18' '  -3.65%  '
Set-TextValue $ws 'D19' '0.06629'
Set-TextValue $ws 'E19' '  -0.33%  '
Set-TextValue $ws 'D20' '17.88'
Set-TextValue $ws 'E20' '  -2.07%  '
Set-TextValue $ws 'D21' '1.002'
Set-TextValue $ws 'E21' '  +0.04%  '
Set-TextValue $ws 'D22' '6.174'
Set-TextValue $ws 'E22' '  -2.40%  '
Set-TextValue $ws 'D23' '28.020.65'
Set-TextValue $ws 'E23' '  -2.31%  '
Set-TextValue $ws 'D24' '11.35'
Set-TextValue $ws 'E24' '  -1.63%  '
Set-TextValue $ws 'D25' '2.315'
Set-TextValue $ws 'E25' '  +0.64%  '
Set-TextValue $ws 'D26' '2.096.83'
Set-TextValue $ws 'E26' '  -1.77%  '
Set-TextValue $ws 'D27' '2.528'
Set-TextValue $ws 'E27' '  -6.09%  '
Set-TextValue $ws 'D28' '157.64'
Set-TextValue $ws 'E28' '  -0.63%  '
Set-TextValue $ws 'E29' '  -2.81%  '
Set-TextValue $ws 'D30' '126.43'
Set-TextValue $ws 'E30' '  -1.82%  '
Set-TextValue $ws 'E31' '  -4.30%  '
Set-TextValue $ws 'D32' '0.1053'
Set-TextValue $ws 'E32' '  -2.98%  '
Set-TextValue $ws 'E33' '  -3.29%  '
Set-TextValue $ws 'D34' '3.594'
Set-TextValue $ws 'E34' '  -1.03%  '
Set-TextValue $ws 'D35' '9.338'
Set-TextValue $ws 'E35' '  -5.81%  '
Set-TextValue $ws 'D36' '0.06537'
Set-TextValue $ws 'E36' '  -3.87%  '
Set-TextValue $ws 'D37' '0.02395'
Set-TextValue $ws 'E37' '  -1.77%  '
Set-TextValue $ws 'B38' 'Algorand'
Set-TextValue $ws 'C38' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D38' '0.2177'
Set-TextValue $ws 'E38' '  -2.75%  '
Set-TextValue $ws 'B39' 'TrustWalletToken'
Set-TextValue $ws 'C39' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws 'D39' '1.288'
Set-TextValue $ws 'E39' '  +8.06%  '
Set-TextValue $ws 'D40' '1.202'
Set-TextValue $ws 'E40' '  -5.06%  '
Set-TextValue $ws 'D41' '0.6388'
Set-TextValue $ws 'E41' '  -1.71%  '
Set-TextValue $ws 'D42' '11.51'
Set-TextValue $ws 'E42' '  -2.53%  '
Set-TextValue $ws 'D43' '4.925'
Set-TextValue $ws 'E43' '  -3.32%  '
Set-TextValue $ws 'D44' '1.001'
Set-TextValue $ws 'E44' '  -0.01%  '
Set-TextValue $ws 'E45' '  -1.94%  '
Set-TextValue $ws 'D46' '0.6010'
Set-TextValue $ws 'E46' '  -1.74%  '
Set-TextValue $ws 'D47' '1.295'
Set-TextValue $ws 'E47' '  +0.63%  '
Set-TextValue $ws 'D48' '3.671'
Set-TextValue $ws 'E48' '  -2.16%  '
Set-TextValue $ws 'D49' '1.985'
Set-TextValue $ws 'E49' '  -2.62%  '
Set-TextValue $ws 'D50' '1.200'
Set-TextValue $ws 'E50' '  -0.72%  '
Set-TextValue $ws 'D51' '120.90'
Set-TextValue $ws 'E51' '  -3.49%  '
